$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# Free up A38's old shared string slot ("DORN (caffe)") so that the new
# note text can reuse it, matching how the shared string table was
# rebuilt when the author re-keyed this section of the table.
$ws.Range("A38").ClearContents()

# --- Row 44 gets a new note in column L ---
$ws.Range("L44").Value = "Does not match what is published"
$ws.Range("L44").WrapText = $true

# --- Row 40 (new row): DORN (GT histogram matching) ---
$ws.Range("A40").Value = "DORN (GT histogram matching)"
$ws.Range("B40").Value = 0.9022
$ws.Range("C40").Value = 0.973
$ws.Range("D40").Value = 0.9905
$ws.Range("F40").Value = 0.4236
$ws.Range("G40").Value = 0.0993
$ws.Range("I40").Value = 0.0417

# --- Row 45 (new row): DenseDepth (GT histogram matching) ---
$ws.Range("A45").Value = "DenseDepth (GT histogram matching)"
$ws.Range("B45").Value = 0.9295
$ws.Range("C45").Value = 0.984
$ws.Range("D45").Value = 0.9954
$ws.Range("F45").Value = 0.338
$ws.Range("G45").Value = 0.0803
$ws.Range("I45").Value = 0.0342

# --- Row 38: now DORN (pytorch), with new data + note ---
$ws.Range("A38").Value = "DORN (pytorch)"
$ws.Range("B38").Value = 0.846
$ws.Range("C38").Value = 0.9542
$ws.Range("D38").Value = 0.9827
$ws.Range("F38").Value = 0.5007
$ws.Range("G38").Value = 0.1197
$ws.Range("I38").Value = 0.0533
$ws.Rows.Item(38).RowHeight = 17

$ws.Range("L38").Value = "DORN evaluation doesn't match but whatever"
$ws.Range("L38").WrapText = $true

# --- Row 39: DORN (median rescaling), fill in data ---
$ws.Range("B39").Value = 0.8683
$ws.Range("C39").Value = 0.9642
$ws.Range("D39").Value = 0.9885
$ws.Range("F39").Value = 0.4838
$ws.Range("G39").Value = 0.1151
$ws.Range("I39").Value = 0.0493

# --- Row 44: DenseDepth (median rescaling), fill in data ---
$ws.Range("B44").Value = 0.8876
$ws.Range("C44").Value = 0.9778
$ws.Range("D44").Value = 0.9952
$ws.Range("F44").Value = 0.4091
$ws.Range("G44").Value = 0.106
$ws.Range("I44").Value = 0.045

# --- number formatting for the new numeric cells to match existing style ---
$ws.Range("B38:D38,F38:G38,I38").NumberFormat = "0.000"
$ws.Range("B39:D39,F39:G39,I39").NumberFormat = "0.000"
$ws.Range("B40:D40,F40:G40,I40").NumberFormat = "0.000"
$ws.Range("B44:D44,F44:G44,I44").NumberFormat = "0.000"
$ws.Range("B45:D45,F45:G45,I45").NumberFormat = "0.000"

# --- sheet view changes ---
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("L38").Select()
